$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.431.70"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "3.766.34"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'593.15"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "'165.79"
$ws.Range("E6").Value = "  -3.11%  "
$ws.Range("D7").Value = "3.764.30"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  -4.72%  "
$ws.Range("D14").Value = "'35.77"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("D15").Value = "4.396.37"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "3.774.71"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "67.345.57"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").Value = "'17.72"
$ws.Range("E18").Value = "  -3.39%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'6.94"
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").Value = "'10.16"
$ws.Range("E21").Value = "  -8.89%  "
$ws.Range("D22").Value = "'455.41"
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("D23").Value = "'0.695"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("E24").Value = "  +4.24%  "
$ws.Range("D25").Value = "'83.07"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").Value = "'11.82"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("D27").Value = "'2.11"
$ws.Range("E27").Value = "  -6.68%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'9.99"
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("D31").Value = "'7.19"
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("D32").Value = "'29.65"
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("D33").Value = "'2.18"
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("D34").Value = "'9.16"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "3.718.71"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "'0.0994"
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("D38").Value = "'3.30"
$ws.Range("E38").Value = "  -7.05%  "
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").Value = "'0.990"
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").Value = "'5.72"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'43.60"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D45").Value = "'0.297"
$ws.Range("E45").Value = "  -4.72%  "
$ws.Range("D46").Value = "'46.75"
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("D48").Value = "'147.35"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("E49").Value = "  -7.49%  "
$ws.Range("D50").Value = "'388.62"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").Value = "2.744.59"
$ws.Range("E51").Value = "  +1.53%  "
